# Apply cell content updates described by the authoritative OOXML diff.
# D (Price) and E (Volume 1h) columns hold text values (e.g. "29.403.72",
# "  +0.06%  ") in the source workbook. Excel would otherwise silently
# coerce numeric-looking strings (e.g. "0.9993") into floating point
# numbers, so force the cell format to Text before writing those values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.403.72'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.26'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.65'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6274'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07495'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2903'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.42'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07747'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.847.85'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.002'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6802'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001038'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.16'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.103.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.184'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.442.31'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '229.41'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.35'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.468'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.78'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1377'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.417'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.55'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06415'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +14.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.390'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.475'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.097'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.067'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.831'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.141'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7001'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.581'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.265.44'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.832'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01831'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.616'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9102'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.007.88'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -18.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.76'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '66.31'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.755'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.091'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1173'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.79%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.007'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3952'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.48%  '
